# remove -zi in name
# Fix the typo "Recognization" -> "Recognition" in the Java class name and
# split the "activity_recognization.xml" label into three runs, correcting
# the middle segment to "_recognition" on the way.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "TextBox 207" shape that holds the RecognitionActivity labels.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText -and $candidate.TextFrame.TextRange.Text -like "*Recognization*") {
        $shape = $candidate
        break
    }
}

$tr = $shape.TextFrame.TextRange

# --- Paragraph 1: RecognizationActivity.java -> RecognitionActivity.java ---
$para1 = $tr.Paragraphs(1)
$run1 = $para1.Runs(1)
$run1.Text = "RecognitionActivity.java"

# --- Paragraph 3: activity_recognization.xml -> activity / _recognition / .xml ---
$para3 = $tr.Paragraphs(3)
$fullText = $para3.Text

$firstLen = "activity".Length
$lastLen = ".xml".Length
$midStart = $firstLen + 1
$midLen = $fullText.Length - $firstLen - $lastLen

$midRange = $para3.Characters($midStart, $midLen)
$midRange.Text = "_recognition"
